$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column (price) values are written as text so strings like
# "1.00", "169.30", "0.0000154" are not silently re-parsed as numbers
# (which would strip trailing zeros / introduce float noise / sci notation).
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '68.966.64'
$ws.Range('E2').Value = '  +3.32%  '
$ws.Range('D3').Value = '3.749.24'
$ws.Range('E3').Value = '  +2.90%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '602.76'
$ws.Range('E5').Value = '  +2.28%  '
$ws.Range('D6').Value = '169.30'
$ws.Range('E6').Value = '  +3.51%  '
$ws.Range('D7').Value = '3.747.36'
$ws.Range('E7').Value = '  +2.93%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  +2.92%  '
$ws.Range('E10').Value = '  +5.67%  '
$ws.Range('D11').Value = '6.35'
$ws.Range('E11').Value = '  +4.39%  '
$ws.Range('D12').Value = '0.463'
$ws.Range('E12').Value = '  +1.47%  '
$ws.Range('D13').Value = '38.26'
$ws.Range('E13').Value = '  +3.26%  '
$ws.Range('E14').Value = '  +4.77%  '
$ws.Range('D15').Value = '4.378.71'
$ws.Range('E15').Value = '  +2.95%  '
$ws.Range('D16').Value = '3.761.88'
$ws.Range('E16').Value = '  +2.89%  '
$ws.Range('D17').Value = '68.939.93'
$ws.Range('E17').Value = '  +3.19%  '
$ws.Range('D18').Value = '7.29'
$ws.Range('E18').Value = '  +3.42%  '
$ws.Range('D19').Value = '0.114'
$ws.Range('E19').Value = '  +0.63%  '
$ws.Range('D20').Value = '17.05'
$ws.Range('E20').Value = '  +1.33%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '10.76'
$ws.Range('E21').Value = '  +19.98%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').Value = '496.89'
$ws.Range('E22').Value = '  +2.48%  '
$ws.Range('D23').Value = '0.728'
$ws.Range('E23').Value = '  +3.13%  '
$ws.Range('D24').Value = '0.0000154'
$ws.Range('E24').Value = '  +13.25%  '
$ws.Range('D25').Value = '85.46'
$ws.Range('E25').Value = '  +0.94%  '
$ws.Range('E26').Value = '  +2.89%  '
$ws.Range('D27').Value = '12.37'
$ws.Range('E27').Value = '  +3.27%  '
$ws.Range('D28').Value = '10.34'
$ws.Range('E28').Value = '  +5.51%  '
$ws.Range('E29').Value = '  +0.69%  '
$ws.Range('D30').Value = '2.54'
$ws.Range('E30').Value = '  +9.18%  '
$ws.Range('D31').Value = '2.98'
$ws.Range('E31').Value = '  +3.41%  '
$ws.Range('D32').Value = '7.93'
$ws.Range('E32').Value = '  +3.97%  '
$ws.Range('D33').Value = '31.90'
$ws.Range('E33').Value = '  +2.11%  '
$ws.Range('D34').Value = '3.893.67'
$ws.Range('E34').Value = '  +2.84%  '
$ws.Range('E35').Value = '  +3.15%  '
$ws.Range('D36').Value = '3.683.50'
$ws.Range('E36').Value = '  +2.89%  '
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('D38').Value = '1.01'
$ws.Range('E38').Value = '  +2.90%  '
$ws.Range('D39').Value = '5.86'
$ws.Range('E39').Value = '  +3.27%  '
$ws.Range('D40').Value = '0.132'
$ws.Range('E40').Value = '  +1.83%  '
$ws.Range('D41').Value = '0.325'
$ws.Range('E41').Value = '  +2.22%  '
$ws.Range('D42').Value = '438.36'
$ws.Range('E42').Value = '  +2.20%  '
$ws.Range('D44').Value = '48.93'
$ws.Range('E44').Value = '  +1.11%  '
$ws.Range('D45').Value = '1.98'
$ws.Range('E45').Value = '  +4.31%  '
$ws.Range('E46').Value = '  +2.88%  '
$ws.Range('D48').Value = '40.63'
$ws.Range('E48').Value = '  +3.03%  '
$ws.Range('D49').Value = '141.25'
$ws.Range('E49').Value = '  -0.59%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = '0.0356'
$ws.Range('E50').Value = '  +4.30%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.781.94'
$ws.Range('E51').Value = '  +1.78%  '
